$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two "Tasks Done" bullet entries for the "شهریور و مهر 99" block ---
$ws.Range("E116").Value = "• Adjusted orientstion widget for the phantom"
$ws.Range("E117").Value = "• Adjusted 2D views for the phantom"

# --- Insert a new activity row for "* 2D/3D Views", pushing the existing rows
#     (Meetings & other / totals / paid / not paid) down by one. ---
$ws.Rows("120:120").Insert()

# Restore border/format for the freshly inserted row's B/C cells (copy from the row above).
$ws.Range("B119:C119").Copy()
$ws.Range("B120:C120").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B120").Value = "* 2D/3D Views"
$ws.Range("C120").Value = 2

# The "Total Hours" row styling (row height + payment formula) ends up shifted onto the
# "Meetings & other" row (121) while the arithmetic total itself lands on row 122.
$ws.Rows(121).RowHeight = 15.6
$ws.Rows(122).EntireRow.AutoFit()

$ws.Range("E122").Copy()
$ws.Range("E121").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E121").Formula = "=C122*40000"
$ws.Range("E122").Clear()

# Shift the "Paid"/"Not Paid" rows: the static Paid value moves up to row 122 and the
# Not-Paid formula moves up to row 123, leaving row 124 without a D value.
$ws.Range("D123").Copy()
$ws.Range("D122").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D122").Value = 0

$ws.Range("D124").Copy()
$ws.Range("D123").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D123").Formula = "=C122-D122"

$ws.Range("D124").Clear()

# --- Update the sheet selection/active cell to reflect the new layout ---
$ws.Range("E124").Select()
